$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.930.19'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '1.815.90'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.97'
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4681'
$ws.Range("E7").Value = '  +0.87%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3691'
$ws.Range("E8").Value = '  -1.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07365'
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8707'
$ws.Range("E10").Value = '  -1.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.40'
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("D12").Value = '1.816.83'
$ws.Range("E12").Value = '  +3.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.382'
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07076'
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.521'
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.63'
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.71'
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").Value = '26.955.71'
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.325'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  -1.75%  '
$ws.Range("D24").Value = '2.061.94'
$ws.Range("E24").Value = '  +3.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.899'
$ws.Range("E25").Value = '  -0.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.62'
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.40'
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.173'
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.337'
$ws.Range("E29").Value = '  +0.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.00'
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08957'
$ws.Range("E31").Value = '  +0.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7700'
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("E33").Value = '  -0.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.507'
$ws.Range("E34").Value = '  +0.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.911'
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.087'
$ws.Range("E37").Value = '  -2.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01963'
$ws.Range("E38").Value = '  +0.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05292'
$ws.Range("E39").Value = '  +1.01%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.272'
$ws.Range("E40").Value = '  +0.72%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.932'
$ws.Range("E41").Value = '  +0.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5316'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.360'
$ws.Range("E43").Value = '  -5.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1658'
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.438'
$ws.Range("E45").Value = '  -2.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4925'
$ws.Range("E46").Value = '  -3.09%  '
$ws.Range("E47").Value = '  +1.11%  '
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.671'
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.81'
$ws.Range("E50").Value = '  -0.75%  '
$ws.Range("E51").Value = '  -0.36%  '
